$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rows that need "x" filled into the (until now empty) columns G through L,
# mirroring the pattern already used on rows 9, 14 and 15.
$rows = @(10, 11, 12, 13, 16, 17, 18, 19)

foreach ($r in $rows) {
    $rangeAddress = "G" + $r + ":L" + $r
    $ws.Range($rangeAddress).Value = "x"
}

# Update the sheet's scroll position and current selection to match the
# saved view state (scrolled down a bit further, selection moved to I18).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I18").Select()
